$d = $word.ActiveDocument

$replacements = @(
    @("34×39=", "87×50="),
    @("62×66=", "33×72="),
    @("51×61=", "13×39="),
    @("63×41=", "26×93="),
    @("76×32=", "73×40="),
    @("67×26=", "42×19="),
    @("63×52=", "50×75="),
    @("28×56=", "48×73="),
    @("37×44=", "93×61="),
    @("60×23=", "68×86="),
    @("87×58=", "72×79="),
    @("67×59=", "13×94="),
    @("86×38=", "66×58="),
    @("90×53=", "47×21="),
    @("48×77=", "59×78="),
    @("72×78=", "29×52="),
    @("43×22=", "75×82="),
    @("45×99=", "56×79="),
    @("44×52=", "41×46="),
    @("59×48=", "18×93="),
    @("33×24=", "27×92="),
    @("50×97=", "24×23="),
    @("81×75=", "43×56="),
    @("54×67=", "48×72="),
    @("79×75=", "68×89=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
